$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 88, shifting existing rows 88:186 down to 89:187.
# This also extends the sheet dimension from A1:R186 to A1:R187 automatically,
# and the new row's D cell inherits the date number format (style) from the
# row above (row 87), matching the style used throughout the D column.
$ws.Rows(88).Insert()

# Populate the newly inserted row 88 with its data.
$ws.Range("A88").Value = 2
$ws.Range("B88").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C88").Value = "Coquimbo"
$ws.Range("D88").Value = 45049
$ws.Range("E88").Value = 4
$ws.Range("F88").Value = 100112024
$ws.Range("G88").Value = "Choclo"
$ws.Range("H88").Value = "Dulce o Americano"
$ws.Range("I88").Value = "Primera"
$ws.Range("J88").Value = 540
$ws.Range("K88").Value = 16000
$ws.Range("L88").Value = 17000
$ws.Range("M88").Value = 16500
$ws.Range("N88").Value = "$/malla 70 unidades"
$ws.Range("O88").Value = "Provincia de Limarí"
$ws.Range("P88").Value = 236
$ws.Range("Q88").Value = 70
$ws.Range("R88").Value = "Hortaliza"
